# Remove the stray trailing "16" that was appended to every Scripture
# reference in column A (rows 2-56), e.g. "Habakkuk 1:216" -> "Habakkuk 1:2".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value2
    if ($value -ne $null -and $value -like "*16") {
        $cell.Value2 = $value.Substring(0, $value.Length - 2)
    }
}
